$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 180; this shifts the existing rows 180-187 down to 181-188
$ws.Rows.Item(180).Insert()

# Populate the newly inserted row 180 with the new weekly record
$ws.Cells.Item(180, 1).Value = 7
$ws.Cells.Item(180, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(180, 3).Value = "Ñuble"
$ws.Cells.Item(180, 4).Value = 44509
$ws.Cells.Item(180, 5).Value = 16
$ws.Cells.Item(180, 6).Value = 100112023
$ws.Cells.Item(180, 7).Value = "Brócoli"
$ws.Cells.Item(180, 8).Value = "Sin especificar"
$ws.Cells.Item(180, 9).Value = "Primera"
$ws.Cells.Item(180, 10).Value = 300
$ws.Cells.Item(180, 11).Value = 650
$ws.Cells.Item(180, 12).Value = 700
$ws.Cells.Item(180, 13).Value = 675
$ws.Cells.Item(180, 14).Value = "$/unidad"
$ws.Cells.Item(180, 15).Value = "Región del Maule"
$ws.Cells.Item(180, 16).Value = 675
$ws.Cells.Item(180, 17).Value = 1
$ws.Cells.Item(180, 18).Value = "Hortaliza"
